$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD1").Value = "time"
$ws.Range("AE1").Value = "plasma_radioactivity"
$ws.Range("AF1").Value = "metabolite_parent_fraction"
$ws.Range("AG1").Value = "whole_blood_radioactivity"

$ws.Range("AD2").Value = 16.200000000000003
$ws.Range("AE2").Value = 32.35522376080818
$ws.Range("AF2").Value = 0.9923034124629081
$ws.Range("AG2").Value = 11.167945849940288
$ws.Range("AD3").Value = 28.2
$ws.Range("AE3").Value = 2507.1051839992842
$ws.Range("AF3").Value = 0.9923034124629081
$ws.Range("AG3").Value = 4027.7881549527583
$ws.Range("AD4").Value = 42
$ws.Range("AE4").Value = 30219.040354280194
$ws.Range("AF4").Value = 0.9923034124629081
$ws.Range("AG4").Value = 36565.981242480346
$ws.Range("AD5").Value = 58.8
$ws.Range("AE5").Value = 49458.133946955189
$ws.Range("AF5").Value = 0.987721971939862
$ws.Range("AG5").Value = 49458.133946955189
$ws.Range("AD6").Value = 73.8
$ws.Range("AE6").Value = 54625.128373263819
$ws.Range("AF6").Value = 0.98277401617497218
$ws.Range("AG6").Value = 77625.006110518632
$ws.Range("AD7").Value = 88.8
$ws.Range("AE7").Value = 59463.809533870204
$ws.Range("AF7").Value = 0.97465638230006379
$ws.Range("AG7").Value = 84130.176446829675
$ws.Range("AD8").Value = 105
$ws.Range("AE8").Value = 64473.636151116443
$ws.Range("AF8").Value = 0.97426988808824744
$ws.Range("AG8").Value = 88701.158697802806
$ws.Range("AD9").Value = 118.8
$ws.Range("AE9").Value = 62047.753161550849
$ws.Range("AF9").Value = 0.9738833938764313
$ws.Range("AG9").Value = 88956.948697178654
$ws.Range("AD10").Value = 133.80000000000001
$ws.Range("AE10").Value = 58211.493736215198
$ws.Range("AF10").Value = 0.97221373888138496
$ws.Range("AG10").Value = 90152.199242522809
$ws.Range("AD11").Value = 148.80000000000001
$ws.Range("AE11").Value = 32321.777348853899
$ws.Range("AF11").Value = 0.96316691140479671
$ws.Range("AG11").Value = 55018.725047026935
$ws.Range("AD12").Value = 181.2
$ws.Range("AE12").Value = 12244.765623047355
$ws.Range("AF12").Value = 0.9758502716844436
$ws.Range("AG12").Value = 35607.826330715769
$ws.Range("AD13").Value = 301.2
$ws.Range("AE13").Value = 7804.7356253636981
$ws.Range("AF13").Value = 0.96181672025723475
$ws.Range("AG13").Value = 26784.480310147766
$ws.Range("AD14").Value = 619.79999999999995
$ws.Range("AE14").Value = 5997.7361438029893
$ws.Range("AF14").Value = 0.92757660167130918
$ws.Range("AG14").Value = 18091.207088138312
$ws.Range("AD15").Value = 915
$ws.Range("AE15").Value = 5549.3370195194557
$ws.Range("AF15").Value = 0.87421987518002875
$ws.Range("AG15").Value = 15843.204817064559
$ws.Range("AD16").Value = 1807.8
$ws.Range("AE16").Value = 4852.2045218107714
$ws.Range("AF16").Value = 0.80636042402826857
$ws.Range("AG16").Value = 12177.168344868838
$ws.Range("AD17").Value = 2710.2000000000003
$ws.Range("AE17").Value = 4121.7383597714706
$ws.Range("AF17").Value = 0.77346278317152117
$ws.Range("AG17").Value = 9460.2027563540141
$ws.Range("AD18").Value = 3610.8
$ws.Range("AE18").Value = 3449.6509856118491
$ws.Range("AF18").Value = 0.70903954802259894
$ws.Range("AG18").Value = 7749.4529025133006
$ws.Range("AD19").Value = 4500
$ws.Range("AE19").Value = 3206.781142108247
$ws.Range("AF19").Value = 0.68246445497630337
$ws.Range("AG19").Value = 7079.0657433029664
$ws.Range("AD20").Value = 5409
$ws.Range("AE20").Value = 3442.7536864959579
$ws.Range("AF20").Value = 0.65989847715736039
$ws.Range("AG20").Value = 5512.0411957519309

$ws.Range("AD1").Font.Size = 11
$ws.Range("AD1").Font.Name = "Calibri"
$ws.Range("AD1").Interior.Color = 5296274

$ws.Range("AD2:AD20").Font.Size = 11
$ws.Range("AD2:AD20").Font.Name = "Calibri"
